# EPBDS-13025: trim the trailing empty array element rendered by the
# sortNull test spreadsheet. The "Value2" result column used to show a
# trailing ", ," (an extra empty element after the sort) - the fix trims
# that so the text ends with a single trailing ", ".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 38-58 in column D hold the rendered "_res_.$Value2$StepN" text.
# All of them end with ", ," except row 56 (the Date step), which has a
# different left-hand side. Trim the duplicated trailing element from
# both distinct strings.
$ws.Range("D38:D58").Value2 = "1, 2, "
$ws.Range("D56").Value2 = "07/12/0080, 07/12/0082, "
